$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 39403.2118634838

$ws.Range("A3:T3").Value = 39375
